$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark task 5 ("Cho phép người choi click chọn các button") as completed
# by filling in the actual start/end date columns (H8, I8) with "20 tháng 10"
$ws.Range("H8").Value = "20 tháng 10"
$ws.Range("I8").Value = "20 tháng 10"

# Update the active selection to I8
$ws.Range("I8").Select()
